{"js": "// Update the date line and all \"A\u00d7B=C\" equation cells in the table,\n// per the diff (answers-of-three-digit_number_multiplied_by_one-digit_number.docx).\n// Each entry is [oldText, newText]; oldText values are unique in the document,\n// so a whole-document search-and-replace for each pair reproduces the diff.\nconst replacements = [\n  [\"2025-01-04 Saturday\", \"2025-01-05 Sunday\"],\n  [\"684\u00d75=3420\", \"813\u00d74=3252\"],\n  [\"773\u00d73=2319\", \"682\u00d79=6138\"],\n  [\"665\u00d77=4655\", \"124\u00d75=620\"],\n  [\"435\u00d74=1740\", \"897\u00d76=5382\"],\n  [\"563\u00d79=5067\", \"157\u00d76=942\"],\n  [\"930\u00d76=5580\", \"717\u00d75=3585\"],\n  [\"815\u00d79=7335\", \"218\u00d79=1962\"],\n  [\"461\u00d72=922\", \"286\u00d73=858\"],\n  [\"255\u00d76=1530\", \"842\u00d72=1684\"],\n  [\"893\u00d75=4465\", \"438\u00d74=1752\"],\n  [\"658\u00d79=5922\", \"877\u00d78=7016\"],\n  [\"177\u00d72=354\", \"952\u00d79=8568\"],\n  [\"457\u00d74=1828\", \"835\u00d76=5010\"],\n  [\"635\u00d73=1905\", \"466\u00d73=1398\"],\n  [\"950\u00d77=6650\", \"435\u00d78=3480\"],\n  [\"720\u00d73=2160\", \"729\u00d75=3645\"],\n  [\"888\u00d78=7104\", \"741\u00d75=3705\"],\n  [\"381\u00d78=3048\", \"592\u00d76=3552\"],\n  [\"674\u00d75=3370\", \"599\u00d79=5391\"],\n  [\"386\u00d77=2702\", \"298\u00d77=2086\"],\n  [\"741\u00d76=4446\", \"710\u00d77=4970\"],\n  [\"824\u00d77=5768\", \"240\u00d72=480\"],\n  [\"420\u00d73=1260\", \"684\u00d76=4104\"],\n  [\"450\u00d79=4050\", \"151\u00d74=604\"],\n  [\"469\u00d78=3752\", \"143\u00d73=429\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all \"A\u00d7B=C\" equation cells in the table,\n# per the diff (answers-of-three-digit_number_multiplied_by_one-digit_number.docx).\n# Each entry is old text / new text; old values are unique in the document,\n# so a document-wide Find/Replace for each pair reproduces the diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-01-04 Saturday\", \"2025-01-05 Sunday\"),\n    @(\"684\u00d75=3420\", \"813\u00d74=3252\"),\n    @(\"773\u00d73=2319\", \"682\u00d79=6138\"),\n    @(\"665\u00d77=4655\", \"124\u00d75=620\"),\n    @(\"435\u00d74=1740\", \"897\u00d76=5382\"),\n    @(\"563\u00d79=5067\", \"157\u00d76=942\"),\n    @(\"930\u00d76=5580\", \"717\u00d75=3585\"),\n    @(\"815\u00d79=7335\", \"218\u00d79=1962\"),\n    @(\"461\u00d72=922\", \"286\u00d73=858\"),\n    @(\"255\u00d76=1530\", \"842\u00d72=1684\"),\n    @(\"893\u00d75=4465\", \"438\u00d74=1752\"),\n    @(\"658\u00d79=5922\", \"877\u00d78=7016\"),\n    @(\"177\u00d72=354\", \"952\u00d79=8568\"),\n    @(\"457\u00d74=1828\", \"835\u00d76=5010\"),\n    @(\"635\u00d73=1905\", \"466\u00d73=1398\"),\n    @(\"950\u00d77=6650\", \"435\u00d78=3480\"),\n    @(\"720\u00d73=2160\", \"729\u00d75=3645\"),\n    @(\"888\u00d78=7104\", \"741\u00d75=3705\"),\n    @(\"381\u00d78=3048\", \"592\u00d76=3552\"),\n    @(\"674\u00d75=3370\", \"599\u00d79=5391\"),\n    @(\"386\u00d77=2702\", \"298\u00d77=2086\"),\n    @(\"741\u00d76=4446\", \"710\u00d77=4970\"),\n    @(\"824\u00d77=5768\", \"240\u00d72=480\"),\n    @(\"420\u00d73=1260\", \"684\u00d76=4104\"),\n    @(\"450\u00d79=4050\", \"151\u00d74=604\"),\n    @(\"469\u00d78=3752\", \"143\u00d73=429\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
